# Add a cross-sectional "Area" calculation next to the existing discharge
# (Q) calculation, plus two small "pickup" cells (J2/K2) that surface the
# running Atotal/Qtotal next to each other for a quick read.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new labels in G1, H1, and the re-used labels in J1/K1 ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Column G: incremental cross-sectional area per segment ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
# G4:G15 share one formula (relative refs step down per row), same pattern
# the sheet already uses for the D6:D19 / E3:E19 shared formulas.
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Column H: running total of the area (only needed once, row 2) ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Columns J/K: quick-glance copies of the totals, row 2 only ---
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- View state: mirror the author's on-screen selection/scroll ---
$ws.Range("G1:K15").Select()
